$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at position 41, shifting existing rows 41-48 down to 42-49
$ws.Rows.Item(41).Insert()

# Fill the new row's content: question label in A41, answer "nein" in B41
$ws.Cells.Item(41, 1).Value = "juenger als 23 oder vor 1940 geboren?"
$ws.Cells.Item(41, 2).Value = "nein"

# Match the style used by row 40 (B column uses an integer number format, style index 10)
$ws.Cells.Item(41, 2).NumberFormat = $ws.Cells.Item(40, 2).NumberFormat

# Update the sheet view (topLeftCell / selection) to match the target state
$ws.Range("A20").Select()
$excel.ActiveWindow.ScrollRow = 10
